# ManageProfile.docx edit:
#  1. In the "Typical Course of Events" table:
#     - Row "Updates necessary details." / "If data in invalid, repeat step 3."
#       -> System response text becomes "Save changes."
#     - The following row (empty actor cell / "Else, new data is valid, save
#       the changes.") is removed entirely.
#  2. In "Alternative Courses": the placeholder line "Line n: <alternative
#     course of action>" is replaced with real text, and the blank paragraph
#     that used to follow it is removed.

$d = $word.ActiveDocument

# --- 1. Table: drop the "Else, new data is valid..." row, reword the row above ---
$table = $d.Tables(1)

for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows($i)
    if ($row.Cells(2).Range.Text -match "Else, new data is valid") {
        $row.Delete()
        break
    }
}

$d.Content.Find.Execute(
    "If data in invalid, repeat step 3.", $true, $false, $false, $false,
    $false, $true, 1, $false, "Save changes.", 2)

# --- 2. Alternative Courses: replace placeholder text ---
$d.Content.Find.Execute(
    "Line n: <alternative course of action>", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "If the new entered details are invalid, the user is asked to re-enter the details.",
    2)

# Remove the now-orphaned blank paragraph that followed the placeholder line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "^If the new entered details") {
        $next = $d.Paragraphs($i + 1)
        if ($next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
        }
        break
    }
}
